# Add 2022-Q3 data
#
# This script:
#  1. Updates the "总计" (summary) sheet with a new leading row for 2022-Q3
#     and shifts all the previously existing rows down by one.
#  2. Inserts a new worksheet named "2022-Q3" right after "总计", containing
#     the per-fund holding detail for that quarter, and renumbers/keeps the
#     rest of the quarterly sheets in their original relative order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the summary sheet ("总计", first sheet)
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make sure row 9 (new last data row) has the same look (bold/bordered index
# column) as the rest of the A column before we populate it.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

$summaryData = @(
    @(0, "2022-Q3", 7, 0.41),
    @(1, "2022-Q2", 9, 0.99),
    @(2, "2022-Q1", 25, 5.12),
    @(3, "2021-Q4", 34, 8.1),
    @(4, "2021-Q3", 3, 0.23),
    @(5, "2021-Q2", 4, 1.52),
    @(6, "2021-Q1", 2, 1.83),
    @(7, "2020-Q4", 2, 1.96)
)

$r = 2
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------------

# Keep a handle on an existing quarterly detail sheet so the new sheet can
# copy its layout/styling (header row + bold bordered index column) exactly.
$styleDonor = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$newSheet.Name = "2022-Q3"

$styleDonor.Range("A1:H8").Copy()
$newSheet.Range("A1:H8").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding detail rows. Columns B-G hold figures that look numeric but
# are stored as text in the source data (e.g. fund codes with leading
# zeros, or percentages kept at fixed precision), so they are written with
# a leading apostrophe to force text entry; columns A and H are genuine
# numbers.
$detailData = @(
    @(0, "200006", "长城消费增值混合",           "5.46", "90.90", "3.91", "0.2135", 5),
    @(1, "519673", "银河康乐股票A",               "1.94", "93.82", "6.10", "0.1183", 3),
    @(2, "002515", "招商丰益灵活配置混合C",       "1.60", "39.36", "1.94", "0.0310", 10),
    @(3, "016018", "银河康乐股票C",               "0.35", "93.82", "6.10", "0.0214", 3),
    @(4, "002514", "招商丰益灵活配置混合A",       "0.60", "39.36", "1.94", "0.0116", 10),
    @(5, "007142", "嘉合稳健增长灵活配置混合C",   "0.41", "85.87", "2.30", "0.0094", 9),
    @(6, "007141", "嘉合稳健增长灵活配置混合A",   "0.28", "85.87", "2.30", "0.0064", 9)
)

$r = 2
foreach ($row in $detailData) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
